$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 309
$ws1.Range("F7").Value = 855
$ws1.Range("F8").Value = 43
$ws1.Range("F9").Value = 509
$ws1.Range("F11").Value = 294
$ws1.Range("F12").Value = 135
$ws1.Range("F14").Value = 231
$ws1.Range("F15").Value = 33
$ws1.Range("F16").Value = 410
$ws1.Range("F17").Value = 6604
$ws1.Range("F20").Value = 20
$ws1.Range("F21").Value = 7547
$ws1.Range("F24").Value = 3387
$ws1.Range("F26").Value = 1481
$ws1.Range("F27").Value = 889
$ws1.Range("F29").Value = 37
$ws1.Range("F31").Value = 70
$ws1.Range("F32").Value = 211
$ws1.Range("F33").Value = 194
$ws1.Range("F34").Value = 1630
$ws1.Range("F35").Value = 7
$ws1.Range("F36").Value = 154
$ws1.Range("F37").Value = 54
$ws1.Range("F39").Value = 1195
$ws1.Range("F40").Value = 1724

# Sheet 3: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 242
$ws3.Range("F4").Value = 75

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 242
$ws4.Range("F5").Value = 75
$ws4.Range("F7").Value = 309
$ws4.Range("F9").Value = 855
$ws4.Range("F10").Value = 43
$ws4.Range("F11").Value = 509
$ws4.Range("F14").Value = 294
$ws4.Range("F15").Value = 135
$ws4.Range("F18").Value = 231
$ws4.Range("F19").Value = 33
$ws4.Range("F20").Value = 410
$ws4.Range("F21").Value = 6604
$ws4.Range("F24").Value = 20
$ws4.Range("F25").Value = 7547
$ws4.Range("F28").Value = 3387
$ws4.Range("F30").Value = 1481
$ws4.Range("F31").Value = 889
$ws4.Range("F33").Value = 37
$ws4.Range("F35").Value = 70
$ws4.Range("F37").Value = 211
$ws4.Range("F38").Value = 194
$ws4.Range("F39").Value = 1630
$ws4.Range("F40").Value = 7
$ws4.Range("F41").Value = 154
$ws4.Range("F42").Value = 54
$ws4.Range("F44").Value = 1195
$ws4.Range("F45").Value = 1724
